$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.78, 0.16, 13.04704260826111, 12.671875),
    @(39.78, 0.16, 17.27826714515686, 15.03125),
    @(20.38, 0.16, 11.11862778663635, 9.453125),
    @(26.1, 0.16, 12.31966423988342, 10.71875),
    @(24.85, 0.36, 23.62781620025635, 22.515625),
    @(39.06, 0.04000000000000001, 14.29457950592041, 12.515625),
    @(37.37, 0.36, 36.94248414039612, 35.90625),
    @(29.36, 0.04000000000000001, 3.250525951385498, 2.125),
    @(24.16, 0.36, 18.95114183425903, 16.75),
    @(46.83, 0.36, 23.58631563186646, 21.640625),
    @(96.31999999999999, 0.36, 40.86225938796997, 38.46875),
    @(94.73999999999999, 0.04000000000000001, 5.92878532409668, 5.3125),
    @(96.31999999999999, 0.36, 47.24303102493286, 45.703125),
    @(59.03, 0.36, 44.20377826690674, 41.921875),
    @(93.62, 0.16, 20.47512221336365, 19.171875),
    @(93.62, 0.16, 30.63012337684631, 28.515625),
    @(96.31999999999999, 0.36, 81.34160351753235, 79.765625),
    @(93.62, 0.16, 20.86597800254822, 19.8125),
    @(39.78, 0.16, 22.68091750144958, 21.125),
    @(39.78, 0.16, 43.15804982185364, 38.28125),
    @(60.22, 0.16, 43.08988738059998, 41.640625),
    @(39.78, 0.16, 14.59283208847046, 14.484375),
    @(39.78, 0.16, 25.07989716529846, 23.78125),
    @(39.78, 0.16, 29.23488664627075, 24.59375),
    @(39.78, 0.16, 20.58482432365417, 18.8125)
)

$startRow = 83
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rowValues[0]
    $ws.Cells.Item($r, 2).Value = $rowValues[1]
    $ws.Cells.Item($r, 3).Value = $rowValues[2]
    $ws.Cells.Item($r, 4).Value = $rowValues[3]
}

